# Apply "Updated the anypoint monitoring roleIds and permissions" change.
# The single "Anypoint Monitoring User" row (row 31) is split into two rows:
#   - row 31: Monitoring Viewer       (same permission marks as before)
#   - row 32: Monitoring Administrator (new row, subset of permission marks)
# Every row from the old row 32 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 32, pushing old rows 32..45 down to 33..46.
$ws.Rows("32:32").Insert()

# Give the freshly inserted row the same thin-border cell style used by every
# other data row (copy formats only from row 31, which already has that style).
$ws.Range("A31:G31").Copy()
$ws.Range("A32:G32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Row 31 ("Anypoint Monitoring" / "Anypoint Monitoring User") becomes
#    "Monitoring Viewer" - all other marks on that row stay the same.
$ws.Range("B31").Value2 = "Monitoring Viewer"

# 3) Populate the new row 32 with the "Monitoring Administrator" permission.
$ws.Range("A32").Value2 = "Anypoint Monitoring"
$ws.Range("B32").Value2 = "Monitoring Administrator"
$ws.Range("C32").Value2 = "X"
$ws.Range("D32").Value2 = "X"
$ws.Range("E32").Value2 = "X"
$ws.Range("F32").Value2 = ""
$ws.Range("G32").Value2 = ""

# Keep the selection roughly where Excel would have left it after this edit.
$ws.Range("E32").Select()

Write-Output "Applied Anypoint Monitoring role split"
